$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab
$ws.Name = "Task Metric"

# --- Row 8: fill in the second iteration's metrics ---
# C8: date text "19/7/2018" typed as text (matches source: stored as shared string, not a date serial)
$ws.Range("C8").NumberFormat = "m/d/yyyy"
$ws.Range("C8").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C8").Value = "19/7/2018"

$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 28
$ws.Range("F8").Formula = "=E8/D8"

$ws.Range("G8").Value = "Estimates are fairly on track and we are fairly on track."
$ws.Range("G8").WrapText = $true

# Row height for row 8 grows because of the wrapped text
$ws.Rows.Item(8).RowHeight = 29

# --- Sheet view: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("G9").Select()
